$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-14 23:48:23'
$ws.Range('H2').Value = '82%'
$ws.Range('E3').Value = '2026-02-14 23:48:25'
$ws.Range('N3').Value = '-9.1 °C 23:29 TU'
$ws.Range('E4').Value = '2026-02-14 23:48:28'
$ws.Range('J4').Value = '998.9 hPa'
$ws.Range('N4').Value = '2.8 °C 23:11 TU'
$ws.Range('O4').Value = '10.3 °C'
$ws.Range('E5').Value = '2026-02-14 23:48:30'
$ws.Range('N5').Value = '-8.6 °C 23:26 TU'
$ws.Range('O5').Value = '-5.5 °C'
$ws.Range('E6').Value = '2026-02-14 23:48:33'
$ws.Range('H6').Value = '70%'
$ws.Range('J6').Value = '998.9 hPa'
$ws.Range('O6').Value = '10.2 °C'
$ws.Range('E7').Value = '2026-02-14 23:48:35'
$ws.Range('J7').Value = '999.1 hPa'
$ws.Range('E8').Value = '2026-02-14 23:48:38'
$ws.Range('J8').Value = '999.0 hPa'
$ws.Range('E9').Value = '2026-02-14 23:48:40'
$ws.Range('O9').Value = '11.5 °C'
$ws.Range('E10').Value = '2026-02-14 23:48:43'
$ws.Range('O10').Value = '9.9 °C'
$ws.Range('E11').Value = '2026-02-14 23:48:45'
$ws.Range('H11').Value = '57%'
$ws.Range('E12').Value = '2026-02-14 23:48:47'
$ws.Range('N12').Value = '9.2 °C 23:28 TU'
$ws.Range('E13').Value = '2026-02-14 23:48:50'
$ws.Range('J13').Value = '1001.9 hPa'
$ws.Range('E14').Value = '2026-02-14 23:48:52'
$ws.Range('E15').Value = '2026-02-14 23:48:55'
$ws.Range('H15').Value = '56%'
$ws.Range('O15').Value = '11.0 °C'
$ws.Range('E16').Value = '2026-02-14 23:48:57'
$ws.Range('H16').Value = '71%'
$ws.Range('O16').Value = '-6.1 °C'
$ws.Range('E17').Value = '2026-02-14 23:48:59'
$ws.Range('E18').Value = '2026-02-14 23:49:02'
$ws.Range('J18').Value = '999.1 hPa'
$ws.Range('O18').Value = '10.6 °C'
$ws.Range('E19').Value = '2026-02-14 23:49:05'
$ws.Range('H19').Value = '72%'
$ws.Range('E20').Value = '2026-02-14 23:49:07'
$ws.Range('H20').Value = '99%'
$ws.Range('E21').Value = '2026-02-14 23:49:10'
$ws.Range('H21').Value = '65%'
$ws.Range('J21').Value = '1001.6 hPa'
$ws.Range('E22').Value = '2026-02-14 23:49:12'
$ws.Range('H22').Value = '81%'
$ws.Range('E23').Value = '2026-02-14 23:49:15'
$ws.Range('I23').Value = '40.9 mm'
$ws.Range('N23').Value = '-9.4 °C 23:03 TU'
$ws.Range('O23').Value = '-6.4 °C'
$ws.Range('E24').Value = '2026-02-14 23:49:17'
$ws.Range('J24').Value = '1003.3 hPa'
$ws.Range('N24').Value = '6.3 °C 23:20 TU'
$ws.Range('E25').Value = '2026-02-14 23:49:19'
$ws.Range('I25').Value = '21.4 mm'
$ws.Range('E26').Value = '2026-02-14 23:49:22'
$ws.Range('E27').Value = '2026-02-14 23:49:24'
$ws.Range('H27').Value = '73%'
$ws.Range('O27').Value = '-3.4 °C'
$ws.Range('E28').Value = '2026-02-14 23:49:27'
$ws.Range('H28').Value = '62%'
$ws.Range('J28').Value = '998.8 hPa'
$ws.Range('O28').Value = '9.3 °C'
$ws.Range('E29').Value = '2026-02-14 23:49:29'
$ws.Range('K29').Value = '8.6 MJ/m2'
$ws.Range('L29').Value = '74.5 km/h - 355º 23:26 TU'
$ws.Range('E30').Value = '2026-02-14 23:49:32'
$ws.Range('J30').Value = '998.7 hPa'
$ws.Range('O30').Value = '11.3 °C'
$ws.Range('E31').Value = '2026-02-14 23:49:34'
$ws.Range('J31').Value = '998.0 hPa'
$ws.Range('N31').Value = '6.7 °C 23:29 TU'
$ws.Range('O31').Value = '9.0 °C'
$ws.Range('E32').Value = '2026-02-14 23:49:37'
$ws.Range('N32').Value = '1.4 °C 23:28 TU'
$ws.Range('E33').Value = '2026-02-14 23:49:39'
$ws.Range('J33').Value = '1001.3 hPa'
$ws.Range('O33').Value = '3.7 °C'
$ws.Range('E34').Value = '2026-02-14 23:49:42'
$ws.Range('H34').Value = '76%'
$ws.Range('E35').Value = '2026-02-14 23:49:44'
$ws.Range('H35').Value = '84%'
$ws.Range('J35').Value = '1005.7 hPa'
$ws.Range('N35').Value = '0.9 °C 23:29 TU'
$ws.Range('O35').Value = '2.5 °C'
$ws.Range('E36').Value = '2026-02-14 23:49:46'
$ws.Range('H36').Value = '55%'
$ws.Range('J36').Value = '999.6 hPa'
$ws.Range('E37').Value = '2026-02-14 23:49:49'
$ws.Range('H37').Value = '60%'
$ws.Range('J37').Value = '999.7 hPa'
$ws.Range('E38').Value = '2026-02-14 23:49:51'
$ws.Range('H38').Value = '78%'
$ws.Range('N38').Value = '5.9 °C 23:27 TU'
$ws.Range('O38').Value = '9.8 °C'
$ws.Range('E39').Value = '2026-02-14 23:49:54'
$ws.Range('H39').Value = '83%'
$ws.Range('E40').Value = '2026-02-14 23:49:56'
$ws.Range('H40').Value = '62%'
$ws.Range('J40').Value = '1002.2 hPa'
$ws.Range('E41').Value = '2026-02-14 23:49:59'
$ws.Range('J41').Value = '1001.0 hPa'
$ws.Range('N41').Value = '9.8 °C 23:29 TU'
$ws.Range('O41').Value = '13.1 °C'
$ws.Range('E42').Value = '2026-02-14 23:50:01'
$ws.Range('H42').Value = '60%'
$ws.Range('E43').Value = '2026-02-14 23:50:03'
$ws.Range('E44').Value = '2026-02-14 23:50:06'
$ws.Range('I44').Value = '38.0 mm'
$ws.Range('N44').Value = '-8.9 °C 23:12 TU'
$ws.Range('E45').Value = '2026-02-14 23:50:09'
$ws.Range('J45').Value = '1008.6 hPa'
$ws.Range('N45').Value = '-1.0 °C 23:26 TU'
$ws.Range('O45').Value = '2.6 °C'
$ws.Range('E46').Value = '2026-02-14 23:50:11'
$ws.Range('J46').Value = '1004.1 hPa'
$ws.Range('N46').Value = '9.0 °C 22:59 TU'
$ws.Range('O46').Value = '11.4 °C'
